$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (shared strings) ---
$ws.Range("A8").Characters(21, 2).Text = "19"
$ws.Range("C9").Characters(27, 8).Text = "5/8/2023"
$ws.Range("C9").Characters(46, 8).Text = "5/14/2023"

# --- Row 30: cells that change from text placeholders to numbers (need NumberFormat) ---
$ws.Range("D30").Value = 2
$ws.Range("D30").NumberFormat = "#,##0"
$ws.Range("E30").Value = -100
$ws.Range("E30").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("G30").Value = 2
$ws.Range("G30").NumberFormat = "#,##0"
$ws.Range("H30").Value = -50
$ws.Range("H30").NumberFormat = "#,##0.0;""-""#,##0.0"

# --- Remaining numeric cell updates ---
# Row 14
$ws.Range("C14").Value = 2
$ws.Range("D14").Value = 4
$ws.Range("E14").Value = -50
$ws.Range("F14").Value = 9
$ws.Range("G14").Value = 12
$ws.Range("H14").Value = -25
$ws.Range("I14").Value = 41
$ws.Range("J14").Value = 54
$ws.Range("K14").Value = -24.074074074074
$ws.Range("L14").Value = -16.326530612244
$ws.Range("M14").Value = 10.810810810810
$ws.Range("N14").Value = -76.162790697674

# Row 15
$ws.Range("C15").Value = 7
$ws.Range("D15").Value = 11
$ws.Range("E15").Value = -36.363636363636
$ws.Range("F15").Value = 41
$ws.Range("G15").Value = 26
$ws.Range("H15").Value = 57.692307692307
$ws.Range("I15").Value = 154
$ws.Range("J15").Value = 144
$ws.Range("K15").Value = 6.944444444444
$ws.Range("L15").Value = 23.2
$ws.Range("M15").Value = 49.514563106796
$ws.Range("N15").Value = -38.888888888888

# Row 16
$ws.Range("C16").Value = 79
$ws.Range("D16").Value = 94
$ws.Range("E16").Value = -15.957446808510
$ws.Range("F16").Value = 321
$ws.Range("G16").Value = 371
$ws.Range("H16").Value = -13.477088948787
$ws.Range("I16").Value = 1564
$ws.Range("J16").Value = 1603
$ws.Range("K16").Value = -2.432938240798
$ws.Range("L16").Value = 38.775510204081
$ws.Range("M16").Value = 5.178211163416
$ws.Range("N16").Value = -73.496017624131

# Row 17
$ws.Range("C17").Value = 156
$ws.Range("D17").Value = 137
$ws.Range("E17").Value = 13.868613138686
$ws.Range("F17").Value = 597
$ws.Range("G17").Value = 539
$ws.Range("H17").Value = 10.760667903525
$ws.Range("I17").Value = 2667
$ws.Range("J17").Value = 2403
$ws.Range("K17").Value = 10.986267166042
$ws.Range("L17").Value = 36.349693251533
$ws.Range("M17").Value = 72.175597159457
$ws.Range("N17").Value = -10.050590219224

# Row 18
$ws.Range("C18").Value = 46
$ws.Range("D18").Value = 55
$ws.Range("E18").Value = -16.363636363636
$ws.Range("F18").Value = 219
$ws.Range("G18").Value = 218
$ws.Range("H18").Value = 0.458715596330
$ws.Range("I18").Value = 1118
$ws.Range("J18").Value = 1093
$ws.Range("K18").Value = 2.287282708142
$ws.Range("L18").Value = 45.572916666666
$ws.Range("M18").Value = 4
$ws.Range("N18").Value = -83.427216128075

# Row 19
$ws.Range("C19").Value = 129
$ws.Range("D19").Value = 113
$ws.Range("E19").Value = 14.159292035398
$ws.Range("F19").Value = 595
$ws.Range("G19").Value = 563
$ws.Range("H19").Value = 5.683836589698
$ws.Range("I19").Value = 2743
$ws.Range("J19").Value = 2839
$ws.Range("K19").Value = -3.381472349418
$ws.Range("L19").Value = 30.681276798475
$ws.Range("M19").Value = 79.515706806282
$ws.Range("N19").Value = 7.274149393820

# Row 20
$ws.Range("C20").Value = 103
$ws.Range("D20").Value = 53
$ws.Range("E20").Value = 94.339622641509
$ws.Range("F20").Value = 453
$ws.Range("G20").Value = 261
$ws.Range("H20").Value = 73.563218390804
$ws.Range("I20").Value = 1977
$ws.Range("J20").Value = 1484
$ws.Range("K20").Value = 33.221024258760
$ws.Range("L20").Value = 119.911012235818
$ws.Range("M20").Value = 178.059071729958
$ws.Range("N20").Value = -65.107659724673

# Row 21
$ws.Range("C21").Value = 522
$ws.Range("D21").Value = 467
$ws.Range("E21").Value = 11.777301927194
$ws.Range("F21").Value = 2235
$ws.Range("G21").Value = 1990
$ws.Range("H21").Value = 12.311557788944
$ws.Range("I21").Value = 10264
$ws.Range("J21").Value = 9620
$ws.Range("K21").Value = 6.694386694386
$ws.Range("L21").Value = 46.148369642602
$ws.Range("M21").Value = 58.151001540832
$ws.Range("N21").Value = -57.689929510697

# Row 22
$ws.Range("C22").Value = 8
$ws.Range("D22").Value = 8
$ws.Range("E22").Value = 0
$ws.Range("F22").Value = 24
$ws.Range("G22").Value = 30
$ws.Range("H22").Value = -20
$ws.Range("I22").Value = 105
$ws.Range("J22").Value = 124
$ws.Range("K22").Value = -15.322580645161
$ws.Range("L22").Value = 28.048780487804
$ws.Range("M22").Value = -13.934426229508

# Row 23
$ws.Range("C23").Value = 33
$ws.Range("D23").Value = 34
$ws.Range("E23").Value = -2.941176470588
$ws.Range("F23").Value = 138
$ws.Range("G23").Value = 128
$ws.Range("H23").Value = 7.8125
$ws.Range("I23").Value = 647
$ws.Range("J23").Value = 543
$ws.Range("K23").Value = 19.152854511970
$ws.Range("L23").Value = 60.945273631840
$ws.Range("M23").Value = 82.768361581920

# Row 24
$ws.Range("C24").Value = 340
$ws.Range("D24").Value = 325
$ws.Range("E24").Value = 4.615384615384
$ws.Range("F24").Value = 1345
$ws.Range("G24").Value = 1398
$ws.Range("H24").Value = -3.791130185979
$ws.Range("I24").Value = 6299
$ws.Range("J24").Value = 6369
$ws.Range("K24").Value = -1.099073637933
$ws.Range("L24").Value = 44.971231300345
$ws.Range("M24").Value = 45.507045507045

# Row 25
$ws.Range("C25").Value = 211
$ws.Range("D25").Value = 197
$ws.Range("E25").Value = 7.106598984771
$ws.Range("F25").Value = 860
$ws.Range("G25").Value = 795
$ws.Range("H25").Value = 8.176100628930
$ws.Range("I25").Value = 3707
$ws.Range("J25").Value = 3533
$ws.Range("K25").Value = 4.924992923860
$ws.Range("L25").Value = 32.867383512544
$ws.Range("M25").Value = -2.447368421052

# Row 26
$ws.Range("C26").Value = 9
$ws.Range("D26").Value = 17
$ws.Range("E26").Value = -47.058823529411
$ws.Range("F26").Value = 61
$ws.Range("G26").Value = 40
$ws.Range("H26").Value = 52.5
$ws.Range("I26").Value = 242
$ws.Range("J26").Value = 251
$ws.Range("K26").Value = -3.585657370517
$ws.Range("L26").Value = 15.238095238095

# Row 27
$ws.Range("C27").Value = 24
$ws.Range("D27").Value = 22
$ws.Range("E27").Value = 9.090909090909
$ws.Range("F27").Value = 81
$ws.Range("G27").Value = 69
$ws.Range("H27").Value = 17.391304347826
$ws.Range("I27").Value = 378
$ws.Range("J27").Value = 325
$ws.Range("K27").Value = 16.307692307692
$ws.Range("L27").Value = 15.596330275229

# Row 28
$ws.Range("D28").Value = 15
$ws.Range("E28").Value = -53.333333333333
$ws.Range("F28").Value = 21
$ws.Range("G28").Value = 40
$ws.Range("H28").Value = -47.5
$ws.Range("I28").Value = 115
$ws.Range("J28").Value = 182
$ws.Range("K28").Value = -36.813186813186
$ws.Range("L28").Value = -30.722891566265
$ws.Range("M28").Value = -19.014084507042
$ws.Range("N28").Value = -73.004694835680

# Row 29
$ws.Range("C29").Value = 6
$ws.Range("D29").Value = 13
$ws.Range("E29").Value = -53.846153846153
$ws.Range("F29").Value = 18
$ws.Range("G29").Value = 31
$ws.Range("H29").Value = -41.935483870967
$ws.Range("I29").Value = 92
$ws.Range("J29").Value = 155
$ws.Range("K29").Value = -40.645161290322
$ws.Range("L29").Value = -38.255033557047
$ws.Range("M29").Value = -22.033898305084
$ws.Range("N29").Value = -75.853018372703

# Row 30
$ws.Range("J30").Value = 17
$ws.Range("K30").Value = -41.176470588235
$ws.Range("L30").Value = -61.538461538461
